$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column Q holds the "archive date" (date the meter was taken out of service).
# Copy the formatting from column P (same style used across the whole column) first,
# then fill in the header/values so the PasteSpecial call can't clobber cell content.
$ws.Range("P1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("P2:P11").Copy()
$ws.Range("Q2:Q11").PasteSpecial(-4122)  # xlPasteFormats

# Set column Q width (Excel quantizes stored width to increments of ~1/6; 23.4 is the
# closest achievable ColumnWidth that serializes to 24.1667, matching the target 24.1719)
$ws.Columns.Item(17).ColumnWidth = 23.4

# Header for the new column
$ws.Range("Q1").Value = "Дата вывода из эксплуатации"

# Only the first data row has an archive date filled in; the rest stay blank
$ws.Range("Q2").Value = "2022-01-25"
